$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 1143 (shifts existing rows 1143+ down by 2)
$ws.Rows("1143:1144").Insert()

# Row 1143 - new record (Rosara, 1a (guarda), Provincia de Talca)
$ws.Range("A1143").Value = 3
$ws.Range("B1143").Value = "Femacal de La Calera"
$ws.Range("C1143").Value = "Coquimbo"
$ws.Range("D1143").Value = 45223
$ws.Range("E1143").Value = 5
$ws.Range("F1143").Value = 100114001
$ws.Range("G1143").Value = "Papa"
$ws.Range("H1143").Value = "Rosara"
$ws.Range("I1143").Value = "1a (guarda)"
$ws.Range("J1143").Value = 120
$ws.Range("K1143").Value = 29000
$ws.Range("L1143").Value = 29000
$ws.Range("M1143").Value = 29000
$ws.Range("N1143").Value = "`$/saco 25 kilos"
$ws.Range("O1143").Value = "Provincia de Talca"
$ws.Range("P1143").Value = 1160
$ws.Range("Q1143").Value = 25
$ws.Range("R1143").Value = "Hortaliza"

# Row 1144 - new record (Rosara, 1a (guarda), Región de O'Higgins)
$ws.Range("A1144").Value = 3
$ws.Range("B1144").Value = "Femacal de La Calera"
$ws.Range("C1144").Value = "Coquimbo"
$ws.Range("D1144").Value = 45223
$ws.Range("E1144").Value = 5
$ws.Range("F1144").Value = 100114001
$ws.Range("G1144").Value = "Papa"
$ws.Range("H1144").Value = "Rosara"
$ws.Range("I1144").Value = "1a (guarda)"
$ws.Range("J1144").Value = 280
$ws.Range("K1144").Value = 28500
$ws.Range("L1144").Value = 29000
$ws.Range("M1144").Value = 28786
$ws.Range("N1144").Value = "`$/saco 25 kilos"
$ws.Range("O1144").Value = "Región de O'Higgins"
$ws.Range("P1144").Value = 1151
$ws.Range("Q1144").Value = 25
$ws.Range("R1144").Value = "Hortaliza"

Write-Output $ws.UsedRange.Address()
